$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "276.63"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "20.89"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "6.228"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06191"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "3.581"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "6.564"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.484"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.8220"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1624"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08221"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.03495"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03103"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09131"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.768"
$ws.Range("E15").Value = "14MCDexMCB"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001632"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04695"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006477"
$ws.Range("E18").Value = "17TigerCashTCHBestin24h"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.006152"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.001069"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.0001501"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.804"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.282"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.01380"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1202"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0002742"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04662"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007053"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1103"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.003522"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01117"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006201"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.7931"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.001592"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00001904"
